# Auto-generated edit script: updates market-board derived profit columns (H:N)
# across multiple sheets, per scheduled market data refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").ClearContents()
$ws.Range("N3").Value = 0
# Row 74
$ws.Range("H74").Value = 5595
$ws.Range("I74").Value = 3992
$ws.Range("J74").Value = 8266.667
$ws.Range("K74").Value = 3992
$ws.Range("L74").Value = 8266.667
$ws.Range("M74").Value = -3056
$ws.Range("N74").Value = -10138.667
# Row 77
$ws.Range("H77").Value = 5595
$ws.Range("I77").Value = 3992
$ws.Range("J77").Value = 8266.667
$ws.Range("K77").Value = 19960
$ws.Range("L77").Value = 41333.335
$ws.Range("M77").Value = -15280
$ws.Range("N77").Value = -50693.335
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("N102").Value = 0

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 982.75
$ws.Range("I2").Value = 507.625
$ws.Range("J2").Value = 1616.25
$ws.Range("K2").Value = 507.625
$ws.Range("L2").Value = 1616.25
$ws.Range("M2").Value = -394.625
$ws.Range("N2").Value = -1842.25
# Row 24
$ws.Range("H24").Value = 20355
$ws.Range("J24").Value = 20355
$ws.Range("L24").Value = 20355
$ws.Range("N24").Value = -21103
# Row 32
$ws.Range("H32").Value = 8275.575
$ws.Range("I32").Value = 7384.0737
$ws.Range("K32").Value = 7384.0737
$ws.Range("M32").Value = -7097.0737
# Row 97
$ws.Range("H97").Value = 740.7917
$ws.Range("I97").Value = 740.7917
$ws.Range("K97").Value = 740.7917
$ws.Range("M97").Value = -244.7917
# Row 100
$ws.Range("H100").Value = 20355
$ws.Range("J100").Value = 20355
$ws.Range("L100").Value = 20355
$ws.Range("N100").Value = -22519
# Row 103
$ws.Range("H103").Value = 25000
$ws.Range("J103").Value = 25000
$ws.Range("L103").Value = 25000
$ws.Range("N103").Value = -27344
# Row 116
$ws.Range("H116").Value = 982.75
$ws.Range("I116").Value = 507.625
$ws.Range("J116").Value = 1616.25
$ws.Range("K116").Value = 507.625
$ws.Range("L116").Value = 1616.25
$ws.Range("M116").Value = 1786.375
$ws.Range("N116").Value = -6204.25
# Row 122
$ws.Range("H122").Value = 863.25
$ws.Range("I122").Value = 754.17645
$ws.Range("K122").Value = 2262.52935
$ws.Range("M122").Value = 187.4706499999998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 982.75
$ws.Range("I3").Value = 507.625
$ws.Range("J3").Value = 1616.25
$ws.Range("K3").Value = 507.625
$ws.Range("L3").Value = 1616.25
$ws.Range("M3").Value = -393.625
$ws.Range("N3").Value = -1844.25
# Row 94
$ws.Range("H94").Value = 817.5862
$ws.Range("I94").Value = 719.65216
$ws.Range("J94").Value = 1193
$ws.Range("K94").Value = 719.65216
$ws.Range("L94").Value = 1193
$ws.Range("M94").Value = -268.65216
$ws.Range("N94").Value = -2095

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2834.9832
$ws.Range("I31").Value = 2752
$ws.Range("J31").Value = 2877.5386
$ws.Range("K31").Value = 2752
$ws.Range("L31").Value = 2877.5386
$ws.Range("M31").Value = -2457
$ws.Range("N31").Value = -3467.5386
# Row 34
$ws.Range("H34").Value = 2834.9832
$ws.Range("I34").Value = 2752
$ws.Range("J34").Value = 2877.5386
$ws.Range("K34").Value = 2752
$ws.Range("L34").Value = 2877.5386
$ws.Range("M34").Value = -2550
$ws.Range("N34").Value = -3281.5386
# Row 132
$ws.Range("H132").Value = 2707.3447
$ws.Range("I132").Value = 1782.1818
$ws.Range("J132").Value = 3272.7222
$ws.Range("K132").Value = 5346.5454
$ws.Range("L132").Value = 9818.1666
$ws.Range("M132").Value = -2816.5454
$ws.Range("N132").Value = -14878.1666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 8733.333
# Row 34
$ws.Range("H34").Value = 1028.826
$ws.Range("I34").Value = 400
$ws.Range("J34").Value = 1057.409
$ws.Range("K34").Value = 1200
$ws.Range("L34").Value = 3172.227
$ws.Range("M34").Value = -1116
$ws.Range("N34").Value = -3340.227
# Row 39
$ws.Range("H39").Value = 2820.7693
$ws.Range("I39").Value = 950
$ws.Range("J39").Value = 3160.9092
$ws.Range("K39").Value = 2850
$ws.Range("L39").Value = 9482.7276
$ws.Range("M39").Value = -2556
$ws.Range("N39").Value = -10070.7276
# Row 41
$ws.Range("H41").Value = 492
$ws.Range("I41").Value = 256
$ws.Range("J41").Value = 1200
$ws.Range("K41").Value = 768
$ws.Range("L41").Value = 3600
$ws.Range("M41").Value = -430
$ws.Range("N41").Value = -4276
# Row 55
$ws.Range("H55").Value = 1780.909
$ws.Range("I55").Value = 545
$ws.Range("J55").Value = 2055.5557
$ws.Range("K55").Value = 1635
$ws.Range("L55").Value = 6166.6671
$ws.Range("M55").Value = -1458
$ws.Range("N55").Value = -6520.6671
# Row 82
$ws.Range("H82").Value = 8000
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
# Row 85
$ws.Range("H85").Value = 8000
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
# Row 127
$ws.Range("H127").Value = 10102220
$ws.Range("J127").Value = 11364910
$ws.Range("L127").Value = 34094730
$ws.Range("N127").Value = -34104650
# Row 131
$ws.Range("H131").Value = 957.9792
$ws.Range("I131").Value = 461.25
$ws.Range("J131").Value = 1206.3438
$ws.Range("K131").Value = 1383.75
$ws.Range("L131").Value = 3619.0314
$ws.Range("M131").Value = 3656.25
$ws.Range("N131").Value = -13699.0314
# Row 132
$ws.Range("H132").Value = 971.625
$ws.Range("I132").Value = 786.2632
$ws.Range("J132").Value = 1676
$ws.Range("K132").Value = 7076.3688
$ws.Range("L132").Value = 15084
$ws.Range("M132").Value = -4546.3688
$ws.Range("N132").Value = -20144

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6010.4644
$ws.Range("I70").Value = 5504.4614
$ws.Range("J70").Value = 6449
$ws.Range("K70").Value = 5504.4614
$ws.Range("L70").Value = 6449
$ws.Range("M70").Value = -5234.4614
$ws.Range("N70").Value = -6989
# Row 73
$ws.Range("H73").Value = 6010.4644
$ws.Range("I73").Value = 5504.4614
$ws.Range("J73").Value = 6449
$ws.Range("K73").Value = 5504.4614
$ws.Range("L73").Value = 6449
$ws.Range("M73").Value = -4568.4614
$ws.Range("N73").Value = -8321
# Row 97
$ws.Range("H97").Value = 1041.037
$ws.Range("I97").Value = 1030.45
$ws.Range("J97").Value = 1071.2858
$ws.Range("K97").Value = 1030.45
$ws.Range("L97").Value = 1071.2858
$ws.Range("M97").Value = -534.45
$ws.Range("N97").Value = -2063.2858

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 56962.05
$ws.Range("I132").Value = 80637.08
$ws.Range("J132").Value = 5666.1665
$ws.Range("K132").Value = 241911.24
$ws.Range("L132").Value = 16998.4995
$ws.Range("M132").Value = -239381.24
$ws.Range("N132").Value = -22058.4995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1760.5686
$ws.Range("I132").Value = 1390.9032
$ws.Range("J132").Value = 2333.55
$ws.Range("K132").Value = 4172.7096
$ws.Range("L132").Value = 7000.650000000001
$ws.Range("M132").Value = -1642.7096
$ws.Range("N132").Value = -12060.65

